$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the last existing data row (A195) down into the new index cells (A196:A199)
$ws.Range("A195").Copy()
$ws.Range("A196:A199").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @(194, "Monza",    "Fiorentina", 2, 1, 0.67, 1.36, 0.88, 1.85, 0, 1, 0.21, 0.49, 0.7,  2),
    @(195, "Atalanta",  "Juventus",  1, 1, 1.43, 1.56, 2.14, 1.42, 0, 0, 0.71, 0.14, 0.84, 2),
    @(196, "Como",      "Milan",     1, 2, 1.43, 1.35, 1.35, 1.47, 0, 0, 0.08, 0.12, 0.2,  3),
    @(197, "Inter",     "Bologna",   2, 2, 1.79, 0.64, 1.71, 0.79, 1, 0, 0.08, 0.15, 0.23, 3)
)

$startRow = 196
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
    $ws.Cells.Item($row, 12).Value = $vals[11]
    $ws.Cells.Item($row, 13).Value = $vals[12]
    $ws.Cells.Item($row, 14).Value = $vals[13]
    $ws.Cells.Item($row, 15).Value = $vals[14]
}
